{"js": "// Update the date stamp and the 25 multiplication problems/answers in the\n// table to the new values. Each \"before\" string occurs exactly once in the\n// document, so a plain text search-and-replace (scoped to the document\n// body) for each pair is unambiguous and safe.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"2024-11-07 Thursday\", \"2024-11-08 Friday\"],\n  [\"840\u00d77=5880\", \"267\u00d73=801\"],\n  [\"993\u00d76=5958\", \"221\u00d75=1105\"],\n  [\"110\u00d75=550\", \"484\u00d74=1936\"],\n  [\"488\u00d77=3416\", \"941\u00d79=8469\"],\n  [\"815\u00d76=4890\", \"924\u00d74=3696\"],\n  [\"215\u00d79=1935\", \"561\u00d77=3927\"],\n  [\"576\u00d73=1728\", \"257\u00d76=1542\"],\n  [\"555\u00d75=2775\", \"159\u00d72=318\"],\n  [\"406\u00d74=1624\", \"848\u00d73=2544\"],\n  [\"254\u00d76=1524\", \"845\u00d76=5070\"],\n  [\"402\u00d74=1608\", \"628\u00d74=2512\"],\n  [\"846\u00d78=6768\", \"589\u00d75=2945\"],\n  [\"687\u00d76=4122\", \"337\u00d78=2696\"],\n  [\"671\u00d77=4697\", \"691\u00d76=4146\"],\n  [\"621\u00d75=3105\", \"957\u00d74=3828\"],\n  [\"838\u00d78=6704\", \"299\u00d79=2691\"],\n  [\"936\u00d77=6552\", \"542\u00d77=3794\"],\n  [\"682\u00d74=2728\", \"825\u00d78=6600\"],\n  [\"987\u00d74=3948\", \"347\u00d79=3123\"],\n  [\"757\u00d78=6056\", \"169\u00d74=676\"],\n  [\"974\u00d78=7792\", \"757\u00d76=4542\"],\n  [\"438\u00d77=3066\", \"811\u00d73=2433\"],\n  [\"331\u00d77=2317\", \"518\u00d75=2590\"],\n  [\"850\u00d75=4250\", \"757\u00d76=4542\"],\n  [\"665\u00d79=5985\", \"338\u00d78=2704\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and the 25 multiplication problems/answers in the\n# table to the new values. Each \"before\" string occurs exactly once in the\n# document, so Find/Replace (ReplaceAll) scoped to the whole document body\n# is unambiguous and safe for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-11-07 Thursday\", \"2024-11-08 Friday\")\n    ,@(\"840\u00d77=5880\", \"267\u00d73=801\")\n    ,@(\"993\u00d76=5958\", \"221\u00d75=1105\")\n    ,@(\"110\u00d75=550\", \"484\u00d74=1936\")\n    ,@(\"488\u00d77=3416\", \"941\u00d79=8469\")\n    ,@(\"815\u00d76=4890\", \"924\u00d74=3696\")\n    ,@(\"215\u00d79=1935\", \"561\u00d77=3927\")\n    ,@(\"576\u00d73=1728\", \"257\u00d76=1542\")\n    ,@(\"555\u00d75=2775\", \"159\u00d72=318\")\n    ,@(\"406\u00d74=1624\", \"848\u00d73=2544\")\n    ,@(\"254\u00d76=1524\", \"845\u00d76=5070\")\n    ,@(\"402\u00d74=1608\", \"628\u00d74=2512\")\n    ,@(\"846\u00d78=6768\", \"589\u00d75=2945\")\n    ,@(\"687\u00d76=4122\", \"337\u00d78=2696\")\n    ,@(\"671\u00d77=4697\", \"691\u00d76=4146\")\n    ,@(\"621\u00d75=3105\", \"957\u00d74=3828\")\n    ,@(\"838\u00d78=6704\", \"299\u00d79=2691\")\n    ,@(\"936\u00d77=6552\", \"542\u00d77=3794\")\n    ,@(\"682\u00d74=2728\", \"825\u00d78=6600\")\n    ,@(\"987\u00d74=3948\", \"347\u00d79=3123\")\n    ,@(\"757\u00d78=6056\", \"169\u00d74=676\")\n    ,@(\"974\u00d78=7792\", \"757\u00d76=4542\")\n    ,@(\"438\u00d77=3066\", \"811\u00d73=2433\")\n    ,@(\"331\u00d77=2317\", \"518\u00d75=2590\")\n    ,@(\"850\u00d75=4250\", \"757\u00d76=4542\")\n    ,@(\"665\u00d79=5985\", \"338\u00d78=2704\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
